# Inserção dos testes de buscas, implementação dos métodos asserts.

$wb = $excel.ActiveWorkbook

# 1. Add new worksheet "Teste2" (new sheets land before the active sheet,
#    so move it after Planilha1 to match the target order / part naming).
#    NOTE: the worksheet COM reference used for the Move() call goes stale
#    afterwards, so re-resolve both sheets by name once the move is done.
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "Teste2"
$wsNew.Move($null, $wb.Worksheets.Item("Planilha1"))

$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Teste2")

# 1b. Mirror Planilha1's page setup (margins/paper) onto the new sheet
$ws2.PageSetup.LeftMargin = 36.850393728
$ws2.PageSetup.RightMargin = 36.850393728
$ws2.PageSetup.TopMargin = 56.692913399999995
$ws2.PageSetup.BottomMargin = 56.692913399999995
$ws2.PageSetup.HeaderMargin = 22.67716464
$ws2.PageSetup.FooterMargin = 22.67716464
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# 2. Update Planilha1!A2 text (GabrielaNomura6 -> GabrielaNomuraa)
$ws1.Range("A2").Value = "GabrielaNomuraa"

# 3. Populate Teste2 with data
$ws2.Range("A1").Value = "NomeBusca"
$ws2.Range("B1").Value = "Elemento"
$ws2.Range("A2").Value = "HP ELITEPAD 1000 G2 TABLET"
$ws2.Range("B2").Value = "BEATS STUDIO 2 OVER-EAR MATTE BLACK HEADPHONES"
$ws2.Range("A3").Value = "Computador"

# 4. Column widths on Teste2 (engine quantizes ColumnWidth to 1/6-char
#    steps, so feed it the inverse of that quantization to land on the
#    closest reachable stored width to the target).
$ws2.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 49.833333333333336

# 5. Column widths on Planilha1 (K, L columns)
$ws1.Columns.Item(11).ColumnWidth = 9.833333333333334
$ws1.Columns.Item(12).ColumnWidth = 10.666666666666666

# 6. Selection on Teste2 ends at A3
$ws2.Range("A3").Select()
